# Generate Report for Handoff
#
# Inserts two new tracked files into the localization-status report:
#   - 013c0ab4-fd77-4262-842b-df4c7e7dc784.md   (new row, right after 717fb483-...)
#   - 6c7127f0-f445-4919-8766-5b113b957ca0.md   (new row, right before .localization-config)
# on all three worksheets (Overview, zh-cn, de-de), both reported as
# "Ready for handoff" / "Include".

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks - row inserts below do not renumber hyperlink
# anchors, so we rebuild them from scratch once all rows are in place.
$ws.Hyperlinks.Delete()

# Row 3 becomes the new 013c0ab4 row; the old rows 3 (064e0f25) and 4
# (.localization-config) each shift down by one.
$ws.Rows.Item(3).Insert()
# Now .localization-config lives on row 5; insert a new row 5 for 6c7127f0
# ahead of it.
$ws.Rows.Item(5).Insert()

$ws.Range("A3").Value = "013c0ab4-fd77-4262-842b-df4c7e7dc784.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A5").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d6346ddcc4d9ec70f50792559cc5da15acd5b4ca/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f75e53aa527251576ebd2dcc7a2681acfc73987c/e2e/013c0ab4-fd77-4262-842b-df4c7e7dc784.md", $missing, $missing, "013c0ab4-fd77-4262-842b-df4c7e7dc784.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/064e0f25-b0eb-4c64-af88-7a51c01369bc.md", $missing, $missing, "064e0f25-b0eb-4c64-af88-7a51c01369bc.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/64b513a3b303b40b6cd81761a6b4a1281cec63aa/e2e/6c7127f0-f445-4919-8766-5b113b957ca0.md", $missing, $missing, "6c7127f0-f445-4919-8766-5b113b957ca0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/d6346ddcc4d9ec70f50792559cc5da15acd5b4ca/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()

$ws.Range("A3").Value = "013c0ab4-fd77-4262-842b-df4c7e7dc784.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "013c0ab4-fd77-4262-842b-df4c7e7dc784.f75e53aa527251576ebd2dcc7a2681acfc73987c.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-08 06:19:05"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A5").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.zh-cn.xlf"
$ws.Range("D5").Value = "2016-03-08 06:19:05"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d6346ddcc4d9ec70f50792559cc5da15acd5b4ca/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea933d980099f45a130d85c8a4265dd2da98c937/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.zh-cn.xlf", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6109b6af6205e6e48e4eaaa6e1af2f0a0d59ed70/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a4f81dd9a9d54c0d0b18d48d27663a0956bbbbed/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.zh-cn.xlf", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f75e53aa527251576ebd2dcc7a2681acfc73987c/e2e/013c0ab4-fd77-4262-842b-df4c7e7dc784.md", $missing, $missing, "013c0ab4-fd77-4262-842b-df4c7e7dc784.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f75e53aa527251576ebd2dcc7a2681acfc73987c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/013c0ab4-fd77-4262-842b-df4c7e7dc784.f75e53aa527251576ebd2dcc7a2681acfc73987c.zh-cn.xlf", $missing, $missing, "013c0ab4-fd77-4262-842b-df4c7e7dc784.f75e53aa527251576ebd2dcc7a2681acfc73987c.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/064e0f25-b0eb-4c64-af88-7a51c01369bc.md", $missing, $missing, "064e0f25-b0eb-4c64-af88-7a51c01369bc.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f80da01613880be026cb58b4d8a6dd21f9762c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.zh-cn.xlf", $missing, $missing, "064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/64b513a3b303b40b6cd81761a6b4a1281cec63aa/e2e/6c7127f0-f445-4919-8766-5b113b957ca0.md", $missing, $missing, "6c7127f0-f445-4919-8766-5b113b957ca0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64b513a3b303b40b6cd81761a6b4a1281cec63aa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.zh-cn.xlf", $missing, $missing, "6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/d6346ddcc4d9ec70f50792559cc5da15acd5b4ca/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()

$ws.Range("A3").Value = "013c0ab4-fd77-4262-842b-df4c7e7dc784.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "013c0ab4-fd77-4262-842b-df4c7e7dc784.f75e53aa527251576ebd2dcc7a2681acfc73987c.de-de.xlf"
$ws.Range("D3").Value = "2016-03-08 06:19:07"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A5").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.de-de.xlf"
$ws.Range("D5").Value = "2016-03-08 06:19:07"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d6346ddcc4d9ec70f50792559cc5da15acd5b4ca/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97fab4c22c91e2bc0520111680769597c6d1e1a7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.de-de.xlf", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/058ce5e7f7f8800c132457e85f341266c8a1948a/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2c9dd680161215b1ff64e05df5e158f02d0e2c82/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.de-de.xlf", $missing, $missing, "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f75e53aa527251576ebd2dcc7a2681acfc73987c/e2e/013c0ab4-fd77-4262-842b-df4c7e7dc784.md", $missing, $missing, "013c0ab4-fd77-4262-842b-df4c7e7dc784.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f75e53aa527251576ebd2dcc7a2681acfc73987c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/013c0ab4-fd77-4262-842b-df4c7e7dc784.f75e53aa527251576ebd2dcc7a2681acfc73987c.de-de.xlf", $missing, $missing, "013c0ab4-fd77-4262-842b-df4c7e7dc784.f75e53aa527251576ebd2dcc7a2681acfc73987c.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/064e0f25-b0eb-4c64-af88-7a51c01369bc.md", $missing, $missing, "064e0f25-b0eb-4c64-af88-7a51c01369bc.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d8955e546467ebaa36cb5c84ac7b35b66c6f7ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.de-de.xlf", $missing, $missing, "064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/64b513a3b303b40b6cd81761a6b4a1281cec63aa/e2e/6c7127f0-f445-4919-8766-5b113b957ca0.md", $missing, $missing, "6c7127f0-f445-4919-8766-5b113b957ca0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64b513a3b303b40b6cd81761a6b4a1281cec63aa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.de-de.xlf", $missing, $missing, "6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/d6346ddcc4d9ec70f50792559cc5da15acd5b4ca/.localization-config", $missing, $missing, ".localization-config") | Out-Null

Write-Host "Generate Report for Handoff: done"
